$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2970
$ws.Range("I64").Value = 2925
$ws.Range("K64").Value = 2925
$ws.Range("M64").Value = -2677
$ws.Range("H67").Value = 2970
$ws.Range("I67").Value = 2925
$ws.Range("K67").Value = 2925
$ws.Range("M67").Value = -2067
$ws.Range("H76").Value = 3300
$ws.Range("I76").Value = 3350
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3350
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -3035
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3300
$ws.Range("I79").Value = 3350
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3350
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2258
$ws.Range("N79").Value = -5184
$ws.Range("H96").Value = 2084.4614
$ws.Range("I96").Value = 1783.3334
$ws.Range("K96").Value = 5350.0002
$ws.Range("M96").Value = -3977.0002
$ws.Range("H112").Value = 5435.162
$ws.Range("J112").Value = 6048.364
$ws.Range("L112").Value = 18145.092
$ws.Range("N112").Value = -20361.092
$ws.Range("H113").Value = 2984.1667
$ws.Range("I113").Value = 2952.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2952.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 301.5
$ws.Range("N113").Value = -9508
$ws.Range("H116").Value = 3116.6667
$ws.Range("I116").Value = 2895.8333
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 2895.8333
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 546.1667000000002
$ws.Range("N116").Value = -10884
$ws.Range("H137").Value = 1723.1578
$ws.Range("I137").Value = 1907
$ws.Range("J137").Value = 1603.2609
$ws.Range("K137").Value = 5721
$ws.Range("L137").Value = 4809.7827
$ws.Range("M137").Value = -3171
$ws.Range("N137").Value = -9909.7827

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1381.3334
$ws.Range("I2").Value = 895
$ws.Range("K2").Value = 895
$ws.Range("M2").Value = -782
$ws.Range("H45").Value = 2260.9443
$ws.Range("I45").Value = 2033.5555
$ws.Range("K45").Value = 2033.5555
$ws.Range("M45").Value = -1656.5555
$ws.Range("H87").Value = 36000
$ws.Range("J87").Value = 36000
$ws.Range("L87").Value = 36000
$ws.Range("N87").Value = -38496
$ws.Range("H90").Value = 36000
$ws.Range("J90").Value = 36000
$ws.Range("L90").Value = 108000
$ws.Range("N90").Value = -120480
$ws.Range("H116").Value = 1381.3334
$ws.Range("I116").Value = 895
$ws.Range("K116").Value = 895
$ws.Range("M116").Value = 1399

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1381.3334
$ws.Range("I3").Value = 895
$ws.Range("K3").Value = 895
$ws.Range("M3").Value = -781
$ws.Range("H105").Value = 2996.6667
$ws.Range("I105").Value = 2996.6667
$ws.Range("K105").Value = 2996.6667
$ws.Range("M105").Value = -1249.6667
$ws.Range("H109").Value = 40249.26
$ws.Range("J109").Value = 40249.26
$ws.Range("L109").Value = 40249.26
$ws.Range("N109").Value = -43023.26

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5220.1694
$ws.Range("I31").Value = 1827.8096
$ws.Range("J31").Value = 6839.25
$ws.Range("K31").Value = 1827.8096
$ws.Range("L31").Value = 6839.25
$ws.Range("M31").Value = -1532.8096
$ws.Range("N31").Value = -7429.25
$ws.Range("H34").Value = 5220.1694
$ws.Range("I34").Value = 1827.8096
$ws.Range("J34").Value = 6839.25
$ws.Range("K34").Value = 1827.8096
$ws.Range("L34").Value = 6839.25
$ws.Range("M34").Value = -1625.8096
$ws.Range("N34").Value = -7243.25
$ws.Range("H62").Value = 4331.467
$ws.Range("I62").Value = 4363.7393
$ws.Range("J62").Value = 3960.3333
$ws.Range("K62").Value = 4363.7393
$ws.Range("L62").Value = 3960.3333
$ws.Range("M62").Value = -3739.7393
$ws.Range("N62").Value = -5208.3333
$ws.Range("H65").Value = 4331.467
$ws.Range("I65").Value = 4363.7393
$ws.Range("J65").Value = 3960.3333
$ws.Range("K65").Value = 21818.6965
$ws.Range("L65").Value = 19801.6665
$ws.Range("M65").Value = -18698.6965
$ws.Range("N65").Value = -26041.6665
$ws.Range("H122").Value = 1829.2142
$ws.Range("I122").Value = 1386
$ws.Range("J122").Value = 2213.3333
$ws.Range("K122").Value = 4158
$ws.Range("L122").Value = 6639.999899999999
$ws.Range("M122").Value = -1708
$ws.Range("N122").Value = -11539.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 668.6129
$ws.Range("I5").Value = 466.72092
$ws.Range("J5").Value = 1125.5264
$ws.Range("K5").Value = 1400.16276
$ws.Range("L5").Value = 3376.5792
$ws.Range("M5").Value = -1288.16276
$ws.Range("N5").Value = -3600.5792
$ws.Range("H12").Value = 73.64
$ws.Range("I12").Value = 33.294117
$ws.Range("K12").Value = 99.882351
$ws.Range("M12").Value = 73.117649
$ws.Range("H98").Value = 167252.83
$ws.Range("J98").Value = 250579.25
$ws.Range("L98").Value = 751737.75
$ws.Range("N98").Value = -754733.75
$ws.Range("H116").Value = 3153.3125
$ws.Range("I116").Value = 2261.4614
$ws.Range("K116").Value = 6784.3842
$ws.Range("M116").Value = -3342.3842
$ws.Range("H119").Value = 333335680
$ws.Range("I119").Value = 500001020
$ws.Range("K119").Value = 1500003060
$ws.Range("M119").Value = -1499998222
$ws.Range("H122").Value = 3605.4443
$ws.Range("J122").Value = 6726.6665
$ws.Range("L122").Value = 60539.9985
$ws.Range("N122").Value = -65439.9985
$ws.Range("H131").Value = 3879.439
$ws.Range("I131").Value = 658.5714
$ws.Range("J131").Value = 4542.5586
$ws.Range("K131").Value = 1975.7142
$ws.Range("L131").Value = 13627.6758
$ws.Range("M131").Value = 3064.2858
$ws.Range("N131").Value = -23707.6758
$ws.Range("H135").Value = 668.6129
$ws.Range("I135").Value = 466.72092
$ws.Range("J135").Value = 1125.5264
$ws.Range("K135").Value = 4200.48828
$ws.Range("L135").Value = 10129.7376
$ws.Range("M135").Value = -1665.48828
$ws.Range("N135").Value = -15199.7376

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8295.064
$ws.Range("I70").Value = 9194.695
$ws.Range("J70").Value = 5708.625
$ws.Range("K70").Value = 9194.695
$ws.Range("L70").Value = 5708.625
$ws.Range("M70").Value = -8924.695
$ws.Range("N70").Value = -6248.625
$ws.Range("H73").Value = 8295.064
$ws.Range("I73").Value = 9194.695
$ws.Range("J73").Value = 5708.625
$ws.Range("K73").Value = 9194.695
$ws.Range("L73").Value = 5708.625
$ws.Range("M73").Value = -8258.695
$ws.Range("N73").Value = -7580.625
$ws.Range("H80").Value = 1881800
$ws.Range("I80").Value = 4501500
$ws.Range("J80").Value = 135333.33
$ws.Range("K80").Value = 4501500
$ws.Range("L80").Value = 135333.33
$ws.Range("M80").Value = -4500502
$ws.Range("N80").Value = -137329.33
$ws.Range("H83").Value = 1881800
$ws.Range("I83").Value = 4501500
$ws.Range("J83").Value = 135333.33
$ws.Range("K83").Value = 22507500
$ws.Range("L83").Value = 676666.6499999999
$ws.Range("M83").Value = -22502508
$ws.Range("N83").Value = -686650.6499999999
$ws.Range("H107").Value = 462.9091
$ws.Range("I107").Value = 465.77777
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 465.77777
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1454.22223
$ws.Range("N107").Value = -4290
$ws.Range("H122").Value = 16582.428
$ws.Range("I122").Value = 26771.75
$ws.Range("J122").Value = 2996.6667
$ws.Range("K122").Value = 80315.25
$ws.Range("L122").Value = 8990.000100000001
$ws.Range("M122").Value = -77865.25
$ws.Range("N122").Value = -13890.0001
$ws.Range("H134").Value = 12875.2
$ws.Range("J134").Value = 12875.2
$ws.Range("L134").Value = 38625.60000000001
$ws.Range("N134").Value = -43695.60000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 98090.5
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()  # remove (was -4002)
$ws.Range("H82").Value = 1898.238
$ws.Range("I82").Value = 1537.75
$ws.Range("J82").Value = 2120.077
$ws.Range("K82").Value = 1537.75
$ws.Range("L82").Value = 2120.077
$ws.Range("M82").Value = -1176.75
$ws.Range("N82").Value = -2842.077
$ws.Range("H84").Value = 98090.5
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()  # remove (was -10008)
$ws.Range("H85").Value = 1898.238
$ws.Range("I85").Value = 1537.75
$ws.Range("J85").Value = 2120.077
$ws.Range("K85").Value = 1537.75
$ws.Range("L85").Value = 2120.077
$ws.Range("M85").Value = -289.75
$ws.Range("N85").Value = -4616.077
$ws.Range("H87").Value = 98189
$ws.Range("J87").Value = 98189
$ws.Range("L87").Value = 98189
$ws.Range("N87").Value = -100435
$ws.Range("H90").Value = 98189
$ws.Range("J90").Value = 98189
$ws.Range("L90").Value = 294567
$ws.Range("N90").Value = -305799
$ws.Range("H122").Value = 5725.95
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 6207
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 18621
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -23521
$ws.Range("H132").Value = 3348.077
$ws.Range("I132").Value = 2888.4443
$ws.Range("J132").Value = 3591.4119
$ws.Range("K132").Value = 8665.332900000001
$ws.Range("L132").Value = 10774.2357
$ws.Range("M132").Value = -6135.332900000001
$ws.Range("N132").Value = -15834.2357
$ws.Range("H136").Value = 2164.4285
$ws.Range("I136").Value = 1963.375
$ws.Range("J136").Value = 2432.5
$ws.Range("K136").Value = 5890.125
$ws.Range("L136").Value = 7297.5
$ws.Range("M136").Value = -3340.125
$ws.Range("N136").Value = -12397.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1117.5
$ws.Range("I107").Value = 1008
$ws.Range("K107").Value = 3024
$ws.Range("M107").Value = -1104
